$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update DISCHARGE_METHOD -> WW_TYPE lookup values (column B) ---
$ws.Range("B6").Value  = "Evaporation"
$ws.Range("B7").Value  = "Irrigation"
$ws.Range("B10").Value = "Wastewater"
$ws.Range("B16").Value = "Surface"
$ws.Range("B17").Value = "Surface"
$ws.Range("B20").Value = "Wastewater"

# --- Add helper / tooltip text in column D for the reuse rows ---
# D17: note for "Reuse: Indirect Potable" (row 17), 10pt Roboto FF202124
$text17  = "Indirect potable reuse: Uses an environmental buffer, such as a lake, river, or a groundwater aquifer, before the water is treated at a drinking water treatment plant. Direct potable reuse: Involves the treatment and distribution of water without an environmental buffer."
$ws.Range("D17").Value = $text17
$ws.Range("D17").Font.Name = "Roboto"
$ws.Range("D17").Font.Color = 2367776
$ws.Range("D17").Font.Size = 10

$run1_17 = "Indirect potable reuse: "
$run2_17 = "Uses an environmental buffer"
$run3_17 = ", such as a lake, river, or a groundwater aquifer, before the water is treated at a drinking water treatment plant. Direct potable reuse: Involves the treatment and distribution of water without an environmental buffer."
$start2_17 = $run1_17.Length + 1
$len2_17   = $run2_17.Length
$start3_17 = $start2_17 + $len2_17
$len3_17   = $run3_17.Length

$bold17 = $ws.Range("D17").Characters($start2_17, $len2_17)
$bold17.Font.Bold = $true

$tail17 = $ws.Range("D17").Characters($start3_17, $len3_17)
$tail17.Font.Name  = "Roboto"
$tail17.Font.Color = 2367776
$tail17.Font.Size  = 10

# D16: note for "Reuse: Other Non-Potable" (row 16), 12pt Roboto FF202124
$text16  = "Non-potable water reuse – Water is captured, treated, and used for non-drinking purposes, such as toilet flushing, clothes washing, and irrigation."
$ws.Range("D16").Value = $text16
$ws.Range("D16").Font.Name = "Roboto"
$ws.Range("D16").Font.Color = 2367776
$ws.Range("D16").Font.Size = 12

$run1_16 = "Non-potable water reuse – Water is "
$run2_16 = "captured, treated, and used for non-"
$run3_16 = "drinking purposes, such as toilet flushing, clothes washing, and irrigation."
$start2_16 = $run1_16.Length + 1
$len2_16   = $run2_16.Length
$start3_16 = $start2_16 + $len2_16
$len3_16   = $run3_16.Length

$bold16 = $ws.Range("D16").Characters($start2_16, $len2_16)
$bold16.Font.Bold = $true

$tail16 = $ws.Range("D16").Characters($start3_16, $len3_16)
$tail16.Font.Name  = "Roboto"
$tail16.Font.Color = 2367776
$tail16.Font.Size  = 12

# --- Match the saved selection state ---
$ws.Range("J25").Select()
